$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.264759063720703
$ws.Range("B1").Value = 2.919392108917236
$ws.Range("C1").Value = 5.554871559143066
$ws.Range("D1").Value = 1.871373891830444
$ws.Range("E1").Value = 1.032833695411682
